$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.275.79"
$ws.Range("E2").Value = '  +4.24%  '
$ws.Range("D3").Value = "'2.367.06"
$ws.Range("E3").Value = '  +1.70%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'309.66"
$ws.Range("E5").Value = '  -0.69%  '
$ws.Range("D6").Value = "'108.23"
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").Value = "'0.616"
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("D10").Value = "'41.03"
$ws.Range("E10").Value = '  +0.90%  '
$ws.Range("D11").Value = "'0.0919"
$ws.Range("E11").Value = '  +0.26%  '
$ws.Range("D12").Value = "'8.46"
$ws.Range("E12").Value = '  -1.20%  '
$ws.Range("E13").Value = '  +1.72%  '
$ws.Range("D14").Value = "'0.979"
$ws.Range("E14").Value = '  -2.57%  '
$ws.Range("D15").Value = "'2.728.78"
$ws.Range("E15").Value = '  +1.78%  '
$ws.Range("D16").Value = "'15.21"
$ws.Range("E16").Value = '  -1.36%  '
$ws.Range("D17").Value = "'2.365.25"
$ws.Range("E17").Value = '  +1.84%  '
$ws.Range("D18").Value = "'45.251.42"
$ws.Range("E18").Value = '  +4.91%  '
$ws.Range("D19").Value = "'14.61"
$ws.Range("E19").Value = '  +12.68%  '
$ws.Range("D20").Value = "'7.30"
$ws.Range("E20").Value = '  -3.15%  '
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("D22").Value = "'73.14"
$ws.Range("E22").Value = '  -1.28%  '
$ws.Range("D23").Value = "'3.48"
$ws.Range("E23").Value = '  -0.62%  '
$ws.Range("D24").Value = "'260.52"
$ws.Range("E24").Value = '  -3.06%  '
$ws.Range("D25").Value = "'2.31"
$ws.Range("E25").Value = '  +2.16%  '
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("D27").Value = "'11.14"
$ws.Range("E27").Value = '  -0.26%  '
$ws.Range("D28").Value = "'7.31"
$ws.Range("E28").Value = '  -4.01%  '
$ws.Range("E29").Value = '  +2.29%  '
$ws.Range("D30").Value = "'0.0968"
$ws.Range("E30").Value = '  +9.17%  '
$ws.Range("D31").Value = "'22.34"
$ws.Range("E31").Value = '  -1.18%  '
$ws.Range("D32").Value = "'37.60"
$ws.Range("E32").Value = '  -3.21%  '
$ws.Range("D33").Value = "'168.96"
$ws.Range("E33").Value = '  +1.13%  '
$ws.Range("E34").Value = '  +6.69%  '
$ws.Range("E35").Value = '  -0.42%  '
$ws.Range("D36").Value = "'0.117"
$ws.Range("E36").Value = '  +4.26%  '
$ws.Range("D37").Value = "'4.78"
$ws.Range("E37").Value = '  +0.96%  '
$ws.Range("D38").Value = "'2.96"
$ws.Range("E38").Value = '  +4.64%  '
$ws.Range("D39").Value = "'3.92"
$ws.Range("E39").Value = '  +3.49%  '
$ws.Range("D40").Value = "'0.0354"
$ws.Range("E40").Value = '  -2.40%  '
$ws.Range("E41").Value = '  +2.25%  '
$ws.Range("D42").Value = "'99.97"
$ws.Range("E42").Value = '  -4.75%  '
$ws.Range("D43").Value = "'0.231"
$ws.Range("E43").Value = '  -1.85%  '
$ws.Range("D44").Value = "'69.47"
$ws.Range("E44").Value = '  -3.08%  '
$ws.Range("D45").Value = "'12.91"
$ws.Range("E45").Value = '  -2.43%  '
$ws.Range("E46").Value = '  +0.30%  '
$ws.Range("D47").Value = "'80.70"
$ws.Range("E47").Value = '  +5.61%  '
$ws.Range("D48").Value = "'112.03"
$ws.Range("E48").Value = '  -1.77%  '
$ws.Range("D49").Value = "'5.53"
$ws.Range("E49").Value = '  +3.88%  '
$ws.Range("D50").Value = "'9.18"
$ws.Range("E50").Value = '  +2.90%  '
$ws.Range("D51").Value = "'1.667.77"
$ws.Range("E51").Value = '  +0.45%  '
